{"js": "const replacements = [\n  [\"2024-08-21 Wednesday\", \"2024-08-22 Thursday\"],\n  [\"727\u00f76=121, 1\", \"571\u00f78=71, 3\"],\n  [\"896\u00f75=179, 1\", \"345\u00f77=49, 2\"],\n  [\"428\u00f75=85, 3\", \"570\u00f73=190, 0\"],\n  [\"228\u00f73=76, 0\", \"107\u00f75=21, 2\"],\n  [\"747\u00f76=124, 3\", \"593\u00f74=148, 1\"],\n  [\"444\u00f78=55, 4\", \"664\u00f77=94, 6\"],\n  [\"156\u00f77=22, 2\", \"954\u00f78=119, 2\"],\n  [\"306\u00f74=76, 2\", \"862\u00f72=431, 0\"],\n  [\"989\u00f76=164, 5\", \"574\u00f73=191, 1\"],\n  [\"646\u00f77=92, 2\", \"448\u00f76=74, 4\"],\n  [\"510\u00f77=72, 6\", \"202\u00f72=101, 0\"],\n  [\"488\u00f77=69, 5\", \"624\u00f75=124, 4\"],\n  [\"782\u00f77=111, 5\", \"165\u00f75=33, 0\"],\n  [\"771\u00f74=192, 3\", \"861\u00f74=215, 1\"],\n  [\"194\u00f75=38, 4\", \"262\u00f75=52, 2\"],\n  [\"586\u00f75=117, 1\", \"607\u00f72=303, 1\"],\n  [\"632\u00f77=90, 2\", \"572\u00f73=190, 2\"],\n  [\"813\u00f79=90, 3\", \"560\u00f76=93, 2\"],\n  [\"942\u00f73=314, 0\", \"337\u00f75=67, 2\"],\n  [\"506\u00f72=253, 0\", \"617\u00f74=154, 1\"],\n  [\"411\u00f75=82, 1\", \"829\u00f75=165, 4\"],\n  [\"980\u00f79=108, 8\", \"924\u00f76=154, 0\"],\n  [\"802\u00f76=133, 4\", \"862\u00f72=431, 0\"],\n  [\"304\u00f72=152, 0\", \"234\u00f74=58, 2\"],\n  [\"892\u00f74=223, 0\", \"681\u00f72=340, 1\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"2024-08-21 Wednesday\"; New=\"2024-08-22 Thursday\"},\n    @{Old=\"727\u00f76=121, 1\"; New=\"571\u00f78=71, 3\"},\n    @{Old=\"896\u00f75=179, 1\"; New=\"345\u00f77=49, 2\"},\n    @{Old=\"428\u00f75=85, 3\"; New=\"570\u00f73=190, 0\"},\n    @{Old=\"228\u00f73=76, 0\"; New=\"107\u00f75=21, 2\"},\n    @{Old=\"747\u00f76=124, 3\"; New=\"593\u00f74=148, 1\"},\n    @{Old=\"444\u00f78=55, 4\"; New=\"664\u00f77=94, 6\"},\n    @{Old=\"156\u00f77=22, 2\"; New=\"954\u00f78=119, 2\"},\n    @{Old=\"306\u00f74=76, 2\"; New=\"862\u00f72=431, 0\"},\n    @{Old=\"989\u00f76=164, 5\"; New=\"574\u00f73=191, 1\"},\n    @{Old=\"646\u00f77=92, 2\"; New=\"448\u00f76=74, 4\"},\n    @{Old=\"510\u00f77=72, 6\"; New=\"202\u00f72=101, 0\"},\n    @{Old=\"488\u00f77=69, 5\"; New=\"624\u00f75=124, 4\"},\n    @{Old=\"782\u00f77=111, 5\"; New=\"165\u00f75=33, 0\"},\n    @{Old=\"771\u00f74=192, 3\"; New=\"861\u00f74=215, 1\"},\n    @{Old=\"194\u00f75=38, 4\"; New=\"262\u00f75=52, 2\"},\n    @{Old=\"586\u00f75=117, 1\"; New=\"607\u00f72=303, 1\"},\n    @{Old=\"632\u00f77=90, 2\"; New=\"572\u00f73=190, 2\"},\n    @{Old=\"813\u00f79=90, 3\"; New=\"560\u00f76=93, 2\"},\n    @{Old=\"942\u00f73=314, 0\"; New=\"337\u00f75=67, 2\"},\n    @{Old=\"506\u00f72=253, 0\"; New=\"617\u00f74=154, 1\"},\n    @{Old=\"411\u00f75=82, 1\"; New=\"829\u00f75=165, 4\"},\n    @{Old=\"980\u00f79=108, 8\"; New=\"924\u00f76=154, 0\"},\n    @{Old=\"802\u00f76=133, 4\"; New=\"862\u00f72=431, 0\"},\n    @{Old=\"304\u00f72=152, 0\"; New=\"234\u00f74=58, 2\"},\n    @{Old=\"892\u00f74=223, 0\"; New=\"681\u00f72=340, 1\"}\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($r.Old, $false, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)\n    if (-not $result) {\n        throw \"Replacement failed for: $($r.Old)\"\n    }\n}\n"}
